# Refresh market-price-derived profit columns (H, I, J, K, L, M, N) across the
# Leve-profit sheets, per the scheduled market data pull.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 66.666664
$ws.Range("I2").Value = 66.666664
$ws.Range("K2").Value = 66.666664
$ws.Range("M2").Value = 46.333336
# Row 19
$ws.Range("H19").Value = 1663
$ws.Range("J19").Value = 1838.4
$ws.Range("L19").Value = 1838.4
$ws.Range("N19").Value = -2188.4
# Row 39
$ws.Range("H39").Value = 361.8889
$ws.Range("I39").Value = 344.625
$ws.Range("K39").Value = 1033.875
$ws.Range("M39").Value = -737.875
# Row 41
$ws.Range("H41").Value = 7334
$ws.Range("I41").Value = 1500
$ws.Range("K41").Value = 1500
$ws.Range("M41").Value = -1060
# Row 69
$ws.Range("H69").Value = 8810.909
$ws.Range("J69").Value = 8810.909
$ws.Range("L69").Value = 26432.727
$ws.Range("N69").Value = -28180.727
# Row 72
$ws.Range("H72").Value = 8810.909
$ws.Range("J72").Value = 8810.909
$ws.Range("L72").Value = 79298.181
$ws.Range("N72").Value = -88034.181
# Row 138
$ws.Range("H138").Value = 5726
$ws.Range("J138").Value = 6035.3887
$ws.Range("L138").Value = 18106.1661
$ws.Range("N138").Value = -28386.1661

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
# Row 122
$ws.Range("H122").Value = 3928.7334
$ws.Range("I122").Value = 2893.1
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 8679.299999999999
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -6229.299999999999
$ws.Range("N122").Value = -22900

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 6271.9287
$ws.Range("I134").Value = 3052.0908
$ws.Range("K134").Value = 9156.2724
$ws.Range("M134").Value = -6621.2724

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 26731.531
$ws.Range("I31").Value = 2300.5454
$ws.Range("J31").Value = 84318.86
$ws.Range("K31").Value = 2300.5454
$ws.Range("L31").Value = 84318.86
$ws.Range("M31").Value = -2005.5454
$ws.Range("N31").Value = -84908.86
# Row 34
$ws.Range("H34").Value = 26731.531
$ws.Range("I34").Value = 2300.5454
$ws.Range("J34").Value = 84318.86
$ws.Range("K34").Value = 2300.5454
$ws.Range("L34").Value = 84318.86
$ws.Range("M34").Value = -2098.5454
$ws.Range("N34").Value = -84722.86
# Row 86
$ws.Range("H86").Value = 7195
$ws.Range("I86").Value = 4154.857
$ws.Range("J86").Value = 9855.125
$ws.Range("K86").Value = 4154.857
$ws.Range("L86").Value = 9855.125
$ws.Range("M86").Value = -3031.857
$ws.Range("N86").Value = -12101.125
# Row 89
$ws.Range("H89").Value = 7195
$ws.Range("I89").Value = 4154.857
$ws.Range("J89").Value = 9855.125
$ws.Range("K89").Value = 20774.285
$ws.Range("L89").Value = 49275.625
$ws.Range("M89").Value = -15158.285
$ws.Range("N89").Value = -60507.625
# Row 99
$ws.Range("H99").Value = 8000
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()
# Row 107
$ws.Range("H107").Value = 1152.1177
$ws.Range("I107").Value = 684
$ws.Range("J107").Value = 2673.5
$ws.Range("K107").Value = 684
$ws.Range("L107").Value = 2673.5
$ws.Range("M107").Value = 1236
$ws.Range("N107").Value = -6513.5
# Row 126
$ws.Range("H126").Value = 8000
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 7879949
$ws.Range("I34").Value = 11112784
$ws.Range("J34").Value = 336666.66
$ws.Range("K34").Value = 33338352
$ws.Range("L34").Value = 1009999.98
$ws.Range("M34").Value = -33338268
$ws.Range("N34").Value = -1010167.98
# Row 39
$ws.Range("H39").Value = 1364.7142
$ws.Range("I39").Value = 1651.5
$ws.Range("J39").Value = 1250
$ws.Range("K39").Value = 4954.5
$ws.Range("L39").Value = 3750
$ws.Range("M39").Value = -4660.5
$ws.Range("N39").Value = -4338
# Row 92
$ws.Range("H92").Value = 399.1
$ws.Range("J92").Value = 399.1
$ws.Range("L92").Value = 1197.3
$ws.Range("N92").Value = -3693.3
# Row 136
$ws.Range("H136").Value = 1677.6666
$ws.Range("I136").Value = 1325.1818
$ws.Range("K136").Value = 3975.5454
$ws.Range("M136").Value = 1124.4546
# Row 137
$ws.Range("H137").Value = 3912.25
$ws.Range("I137").Value = 1240
$ws.Range("J137").Value = 8366
$ws.Range("K137").Value = 3720
$ws.Range("L137").Value = 25098
$ws.Range("M137").Value = 1380
$ws.Range("N137").Value = -35298
# Row 138
$ws.Range("H138").Value = 4322.4165
$ws.Range("I138").Value = 967.1429000000001
$ws.Range("K138").Value = 2901.4287
$ws.Range("M138").Value = 2238.5713
# Row 139
$ws.Range("H139").Value = 4508.0557
$ws.Range("I139").Value = 1983.3334
$ws.Range("J139").Value = 9557.5
$ws.Range("K139").Value = 5950.0002
$ws.Range("L139").Value = 28672.5
$ws.Range("M139").Value = -810.0002000000004
$ws.Range("N139").Value = -38952.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 6945.857
$ws.Range("I70").Value = 6562.4
$ws.Range("K70").Value = 6562.4
$ws.Range("M70").Value = -6292.4
# Row 73
$ws.Range("H73").Value = 6945.857
$ws.Range("I73").Value = 6562.4
$ws.Range("K73").Value = 6562.4
$ws.Range("M73").Value = -5626.4
# Row 132
$ws.Range("H132").Value = 4189.967
$ws.Range("I132").Value = 1809.25
$ws.Range("K132").Value = 5427.75
$ws.Range("M132").Value = -2897.75
# Row 141
$ws.Range("H141").Value = 219000
$ws.Range("J141").Value = 219000
$ws.Range("L141").Value = 219000
$ws.Range("N141").Value = -229360

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 4327
$ws.Range("I7").Value = 3947.25
$ws.Range("J7").Value = 4833.3335
$ws.Range("K7").Value = 3947.25
$ws.Range("L7").Value = 4833.3335
$ws.Range("M7").Value = -3835.25
$ws.Range("N7").Value = -5057.3335
# Row 106
$ws.Range("H106").Value = 9913
$ws.Range("J106").Value = 9913
$ws.Range("L106").Value = 9913
$ws.Range("N106").Value = -12437
# Row 122
$ws.Range("H122").Value = 6355.6113
$ws.Range("I122").Value = 4142.769
$ws.Range("K122").Value = 12428.307
$ws.Range("M122").Value = -9978.307000000001
# Row 126
$ws.Range("H126").Value = 4327
$ws.Range("I126").Value = 3947.25
$ws.Range("J126").Value = 4833.3335
$ws.Range("K126").Value = 11841.75
$ws.Range("L126").Value = 14500.0005
$ws.Range("M126").Value = -9371.75
$ws.Range("N126").Value = -19440.0005

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
# Row 109
$ws.Range("H109").Value = 67989.75
$ws.Range("J109").Value = 67989.75
$ws.Range("L109").Value = 67989.75
$ws.Range("N109").Value = -70763.75
# Row 113
$ws.Range("H113").Value = 1583.9375
$ws.Range("I113").Value = 1754.9
$ws.Range("J113").Value = 1299
$ws.Range("K113").Value = 5264.700000000001
$ws.Range("L113").Value = 3897
$ws.Range("M113").Value = -3094.700000000001
$ws.Range("N113").Value = -8237
# Row 126
$ws.Range("H126").Value = 4507.7856
$ws.Range("J126").Value = 4672.1113
$ws.Range("L126").Value = 14016.3339
$ws.Range("N126").Value = -18956.3339
